# Structured Phrase Format Descriptions.pptx - apply commit diff
#
# Summary of changes:
#  - Slide 5 (Enrollment): minor wording fix "dates" -> "date"
#  - Slide 6 (Unenrolled): rename the "#clinical trial" / "#enrolled on" /
#    "#ended on" hashtags to the new "#unenrolled" / "#on" enrollment
#    hashtags, reflow several shapes/connectors that shift left as a
#    result of the shorter tags, and delete the now-redundant
#    "#ended on" / "Either ... or ..." explanation shapes (and their
#    connector) that covered the removed phrase variant.

$p = $ppt.ActivePresentation

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Slide 5 - "Enrollment:" -- shape id 16, last run "dates" -> "date"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5shp16 = Get-ShapeById $s5 16
$tr = $s5shp16.TextFrame.TextRange
$tr.Runs(3).Text = " date can be specified in any order.   "

# ---------------------------------------------------------------------
# Slide 6 - "Unenrolled:"
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)

# --- shape id 5: "Patient in #clinical trial #<trial name> #enrolled on #01/01/2017."
#     -> "Patient #unenrolled from #<trial name> #on #01/01/2017."
$shp5 = Get-ShapeById $s6 5
$tr = $shp5.TextFrame.TextRange
$tr.Runs(1).Text = "Patient "
$tr.Runs(2).Text = "#unenrolled"
$tr.Runs(3).Text = " from "
$tr.Runs(6).Text = "#on"

# --- shape id 11: "1. Must start with the #clinical trial tag."
#     -> "1. Must start with the #unenrolled tag."; also moves left.
$shp11 = Get-ShapeById $s6 11
$tr = $shp11.TextFrame.TextRange
$tr.Runs(2).Text = "#unenrolled"
$shp11.Left = 2919957 / 12700

# --- shape id 12: "2. Specify one clinical trial name" (text unchanged, moves right)
$shp12 = Get-ShapeById $s6 12
$shp12.Left = 4841947 / 12700

# --- shape id 13: "3. The #enrolled on date is the date the patient was
#     enrolled in the specified clinical trial." ->
#     "3. The # on date is the date the patient was unenrolled from the
#     specified clinical trial."; also moves right.
$shp13 = Get-ShapeById $s6 13
$tr = $shp13.TextFrame.TextRange
$tr.Runs(2).Text = "# on"
$tr.Runs(3).Text = " date is the date the patient was unenrolled from the specified clinical trial."
$shp13.Left = 6281716 / 12700

# --- shape id 15 (German example): merge "The "/"patient"/" in " runs into
#     a single "The patient " run, then
#     "#clinical trial" -> "#unenrolled", "#ended on" -> "#on".
$shp15 = Get-ShapeById $s6 15
$tr = $shp15.TextFrame.TextRange
$tr.Characters(1, 15).Text = "The patient "
$tr.Runs(2).Text = "#unenrolled"
$tr.Runs(6).Text = "#on"

# --- shape id 16: "Trial name, #enrolled on, and #ended on dates can be
#     specified in any order.   " ->
#     "Trial name, and #on date can be specified in any order.   "
#     also moves right.
$shp16 = Get-ShapeById $s6 16
$tr = $shp16.TextFrame.TextRange
$tr.Characters(40, 40).Text = " date can be specified in any order.   "
$tr.Characters(31, 9).Text = ""
$tr.Characters(25, 6).Text = ""
$tr.Runs(2).Text = "#on"
$tr.Runs(1).Text = "Trial name, and "
$shp16.Left = 3088479 / 12700

# --- connectors that shift along with the boxes above
$conn24 = Get-ShapeById $s6 24
$conn24.Left = 4883035 / 12700

$conn27 = Get-ShapeById $s6 27
$conn27.Left = 6381257 / 12700

$conn29 = Get-ShapeById $s6 29
$conn29.Left = 3007157 / 12700

# --- delete the shapes describing the now-removed "#ended on" / "Either"
#     explanatory boxes and their connecting lines; the remaining
#     "Straight Connector 25" (id 26) stays but moves to a new position.
$toDelete = 21, 25, 19, 20
foreach ($id in $toDelete) {
    $shp = Get-ShapeById $s6 $id
    if ($shp -ne $null) {
        $shp.Delete()
    }
}

$conn26 = Get-ShapeById $s6 26
$conn26.Left = 3179891 / 12700
$conn26.Top = 4477812 / 12700
